$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 10:24 AM"

# --- Top Gainers sheet ---
$wsGain = $wb.Worksheets.Item("Top Gainers")

# Block 1: rows 12-15 re-rank, BUTTERFLY jumps to the top with refreshed figures
$wsGain.Range("B12").Value = "BUTTERFLY"
$wsGain.Range("C12").Value = 7.3915
$wsGain.Range("D12").Value = 10.2954
$wsGain.Range("E12").Value = 12.8889

$wsGain.Range("B13").Value = "VAIBHAVGBL"
$wsGain.Range("C13").Value = 7.3858
$wsGain.Range("D13").Value = 7.9211
$wsGain.Range("E13").Value = 14.3593

$wsGain.Range("B14").Value = "ABREL"
$wsGain.Range("C14").Value = 7.2702
$wsGain.Range("D14").Value = 8.048
$wsGain.Range("E14").Value = 7.6003

$wsGain.Range("B15").Value = "GRAPHITE"
$wsGain.Range("C15").Value = 6.8186
$wsGain.Range("D15").Value = 13.3166
$wsGain.Range("E15").Value = 13.5206

# Block 2: rows 71-74 re-rank, SHK drops to the bottom with refreshed figures
$wsGain.Range("B71").Value = "BHARTIHEXA"
$wsGain.Range("C71").Value = 3.6718
$wsGain.Range("D71").Value = 7.0877
$wsGain.Range("E71").Value = 15.3332

$wsGain.Range("B72").Value = "HLEGLAS"
$wsGain.Range("C72").Value = 3.659
$wsGain.Range("D72").Value = 8.115500000000001
$wsGain.Range("E72").Value = 27.1239

$wsGain.Range("B73").Value = "RHIM"
$wsGain.Range("C73").Value = 3.6544
$wsGain.Range("D73").Value = 3.2276
$wsGain.Range("E73").Value = 5.1826

$wsGain.Range("B74").Value = "SHK"
$wsGain.Range("C74").Value = 3.6347
$wsGain.Range("D74").Value = 2.388
$wsGain.Range("E74").Value = -1.932

# --- Top Losers sheet ---
$wsLose = $wb.Worksheets.Item("Top Losers")

# Rows 70-76 re-rank, DHARMAJ drops out and CSBBANK enters at the bottom
$wsLose.Range("B70").Value = "FINOPB"
$wsLose.Range("C70").Value = -2.3673
$wsLose.Range("D70").Value = -6.2696
$wsLose.Range("E70").Value = 11.1938

$wsLose.Range("B71").Value = "UNIMECH"
$wsLose.Range("C71").Value = -2.353
$wsLose.Range("D71").Value = -1.1572
$wsLose.Range("E71").Value = 0

$wsLose.Range("B72").Value = "FCL"
$wsLose.Range("C72").Value = -2.3453
$wsLose.Range("D72").Value = -2.616
$wsLose.Range("E72").Value = -0.02

$wsLose.Range("B73").Value = "DEEDEV"
$wsLose.Range("C73").Value = -2.3136
$wsLose.Range("D73").Value = -6.6339
$wsLose.Range("E73").Value = -7.4039

$wsLose.Range("B74").Value = "WEALTH"
$wsLose.Range("C74").Value = -2.3047
$wsLose.Range("D74").Value = -3.8606
$wsLose.Range("E74").Value = -2.8234

$wsLose.Range("B75").Value = "RATNAMANI"
$wsLose.Range("C75").Value = -2.2788
$wsLose.Range("D75").Value = -0.4626
$wsLose.Range("E75").Value = 0.8712

$wsLose.Range("B76").Value = "CSBBANK"
$wsLose.Range("C76").Value = -2.2695
$wsLose.Range("D76").Value = 2.3137
$wsLose.Range("E76").Value = 10.6999

# --- 1 Month Performance sheet ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")

# Rows 22-25 re-rank, ONMOBILE jumps to the top with a refreshed figure
$wsPerf.Range("B22").Value = "ONMOBILE"
$wsPerf.Range("C22").Value = 35.4702

$wsPerf.Range("B23").Value = "SHAREINDIA"
$wsPerf.Range("C23").Value = 35.3207

$wsPerf.Range("B24").Value = "SOUTHBANK"
$wsPerf.Range("C24").Value = 35.2819

$wsPerf.Range("B25").Value = "TVSELECT"
$wsPerf.Range("C25").Value = 35.1983
